$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# --- ammo "wgt" (weight) ratio scratch table, rows 17-28 ---
$ws.Range("T17").Value = "wgt"

# Row 18 typed individually (no fill yet -> plain, non-shared formulas)
$ws.Range("T18").Value = 10
$ws.Range("U18").Value = 45
$ws.Range("V18").Formula = "=T18+U18"
$ws.Range("W18").Formula = "=(T18/V18)*2"
$ws.Range("X18").Formula = "=T18/U18"

# Rows 19-22 entered, then the formulas filled down together (one shared group)
$ws.Range("T19").Value = 15
$ws.Range("T20").Value = 30
$ws.Range("T21").Value = 45
$ws.Range("T22").Value = 55
$ws.Range("U19:U22").Value = 45
$ws.Range("V19:V22").Formula = "=T19+U19"
$ws.Range("W19:W22").Formula = "=(T19/V19)*2"
$ws.Range("X19:X22").Formula = "=T19/U19"

# Rows 23-28 entered, then the formulas filled down together (second shared group)
$ws.Range("T23").Value = 75
$ws.Range("T24").Value = 90
$ws.Range("T25").Value = 125
$ws.Range("T26").Value = 300
$ws.Range("T27").Value = 450
$ws.Range("T28").Value = 55
$ws.Range("U23:U28").Value = 45
$ws.Range("V23:V28").Formula = "=T23+U23"
$ws.Range("W23:W28").Formula = "=(T23/V23)*2"
$ws.Range("X23:X28").Formula = "=T23/U23"

# --- leave Sheet4 as the active tab/selection, matching the saved view state ---
$ws.Activate() | Out-Null
$ws.Range("T28").Select() | Out-Null
